$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: rename "Dose" to "Dose (mg)", add "Location" and "Lot" headers
$ws.Range("C1").Value = "Dose (mg)"
$ws.Range("D1").Value = "Location"
$ws.Range("E1").Value = "Lot"

# Existing rows: dose column becomes numeric (mg value only), add Location
$ws.Range("C2").Value = 2.5
$ws.Range("D2").Value = "R Thigh"

$ws.Range("C3").Value = 2.5
$ws.Range("D3").Value = "L Stomach"

# New rows of recent injection data
$ws.Range("A4").Value = 45566.5
$ws.Range("A4").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
$ws.Range("B4").Value = "Zepbound (tirzepatide)"
$ws.Range("C4").Value = 2.5
$ws.Range("D4").Value = "R Stomach"

$ws.Range("A5").Value = 45573.5
$ws.Range("A5").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
$ws.Range("B5").Value = "Zepbound (tirzepatide)"
$ws.Range("C5").Value = 2.5
$ws.Range("D5").Value = "B Stomach"

# Set column D width (closest achievable value to source's 10.140625 "bestFit"
# width; the sandboxed engine quantizes ColumnWidth to internal pixel steps)
$ws.Columns.Item(4).ColumnWidth = 9.25

# Select E5 like the author left off
$null = $ws.Range("E5").Select()
